# Table17.xlsx — wrangling fixes for the "Nativity of fishermen" table.
#
# 1. Fill in the previously-blank "Total" sub-rows for Santa Barbara, Los
#    Angeles and San Diego with the missing "United States" label (and a
#    stray lower-case typo for San Diego, matching the source scan).
# 2. Correct the "Los Angelos" misspelling to "Los Angeles" across all of
#    its rows.
# 3. Clean up the garbled "Jugoslavia —" / "Italy  —---" entries to plain
#    "Jugoslavia" / "Italy".
# 4. Move the active selection to B31 (no more pinned top-left scroll cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Santa Barbara block (rows 21-23): fill missing "United States" total row
$ws.Range("B21").Value = "United States "

# --- Los Angeles block (rows 24-29): fix spelling + fill missing total row +
#     clean up garbled country names
$ws.Range("A24").Value = "Los Angeles "
$ws.Range("B24").Value = "United States "
$ws.Range("A25").Value = "Los Angeles "
$ws.Range("B25").Value = "Jugoslavia"
$ws.Range("A26").Value = "Los Angeles "
$ws.Range("B26").Value = "Italy"
$ws.Range("A27").Value = "Los Angeles "
$ws.Range("A28").Value = "Los Angeles "
$ws.Range("A29").Value = "Los Angeles "

# --- San Diego block (rows 30-35): fill missing total row (typo preserved
#     from the original scan)
$ws.Range("B30").Value = "united States "

# --- Move selection/cursor to where editing left off
$ws.Range("B31").Select() | Out-Null
